# Populate the "MyTickets" Blank-Data report: for the current user, list every
# ticket that has a blank "Serial Number" and/or a blank "Issue Type" field.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MyTickets")

# Start from a clean slate so a re-run doesn't leave stale rows behind.
$ws.Range("A2:C100").ClearContents()

$username = "dmech+greco@barcodesinc.com"

# Ticket id -> which column(s) were found blank for that ticket
$tickets = @(
    @{ Id = "SC1659465"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1655543"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1648899"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1625129"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1625041"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1619781"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1614845"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1614740"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1607016"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1589419"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "ID1858-11"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-10"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-9"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-7"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-8"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-6"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-4"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-5"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-3"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-1"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1858-2"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1859-7"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1859-4"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1859-5"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1859-6"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "ID1859-2"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1537274"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1530937"; BlankSerialNumber = $true; BlankIssueType = $true },
    @{ Id = "SC1520339"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1482696"; BlankSerialNumber = $false; BlankIssueType = $true },
    @{ Id = "SC1482694"; BlankSerialNumber = $false; BlankIssueType = $true }
)

$row = 2
$isFirstRow = $true
foreach ($ticket in $tickets) {
    if ($ticket.BlankSerialNumber) {
        if ($isFirstRow) {
            $ws.Range("A$row").Value = $username
            $isFirstRow = $false
        }
        $ws.Range("B$row").Value = $ticket.Id
        $ws.Range("C$row").Value = "Serial Number"
        $row++
    }
    if ($ticket.BlankIssueType) {
        if ($isFirstRow) {
            $ws.Range("A$row").Value = $username
            $isFirstRow = $false
        }
        $ws.Range("B$row").Value = $ticket.Id
        $ws.Range("C$row").Value = "Issue Type"
        $row++
    }
}

$ws.Range("B$row").Value = "There are no Blank data"
